$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values
$ws.Range("B2").Value = 390
$ws.Range("B7").Value = 405

# Add new row 16 with the new note text, reusing the style already applied
# to the other note rows (e.g. A15) so no new style entries are created.
$ws.Range("A16").Value = "08.03.2025 - Otistics Cj+Tobias karşılığında NSY'ye 4 dolar vermiştir. (390-405)"
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Update selection to A14 as per the saved workbook view
$ws.Range("A14").Select()
